# Transaction download fix with cloudinary
# Adds two new transaction rows (58 and 59) to the Transactions sheet,
# mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58 -----------------------------------------------------------
# Seed row 58 by copying row 57 (value + style + number formats), then
# overwrite the cells with the new transaction's data.
$ws.Range("A57:I57").Copy($ws.Range("A58:I58"))

$ws.Range("A58").Value = "65d8144d9dcadcc6f33c77bd"
$ws.Range("C58").Value = "mamaraffi"
$ws.Range("D58").Value = "Web"
$ws.Range("E58").Value = "Peci Brokat Muslim - Base - L`nPeci Brokat Muslim - Base - XL"
$ws.Range("F58").Value = "Rp. 0"
$ws.Range("G58").Value = "Rp. 58,900"
$ws.Range("H58").Value = "Rp. 58,900"
$ws.Range("I58").Value = "Rp. 0"

# The Date column auto-detects as a date serial when assigned a plain
# date-shaped string, so force it to text first, then restore the
# General-style formatting that the rest of the row uses.
$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = "2/23/2024"
$ws.Range("A58").Copy()
$ws.Range("B58").PasteSpecial(-4122)

# Wrapped product text inflates the row height automatically; AutoFit
# brings it back down to the sheet's normal (default) row height so no
# stray explicit height sticks around.
$ws.Rows(58).AutoFit()

# --- Row 59 -----------------------------------------------------------
$ws.Range("A57:I57").Copy($ws.Range("A59:I59"))

$ws.Range("A59").Value = "65ddfaa404ada948c61fa88e"
$ws.Range("C59").Value = "mamaraffi"
$ws.Range("D59").Value = "Web"
$ws.Range("E59").Value = "Gaun Pesta Elegan - Satin Coklat - One SIze"
$ws.Range("F59").Value = "Rp. 0"
$ws.Range("G59").Value = "Rp. 297,000"
$ws.Range("H59").Value = "Rp. 297,000"
$ws.Range("I59").Value = "Rp. 0"

$ws.Range("B59").NumberFormat = "@"
$ws.Range("B59").Value = "2/27/2024"
$ws.Range("A59").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Rows(59).AutoFit()
